$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.025494267911823
$ws.Range("D2").Value = 1.029622293557509
$ws.Range("E2").Value = 1.025842403853385
$ws.Range("F2").Value = 1.023994663069089
$ws.Range("I2").Value = 1.032300846081392
$ws.Range("J2").Value = 1.030662975878256
$ws.Range("K2").Value = 1.032435670204791
$ws.Range("L2").Value = 1.028666791903852
$ws.Range("M2").Value = 1.026824466054024
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.026368869764338
$ws.Range("D3").Value = 1.03025259309959
$ws.Range("E3").Value = 1.026582943201328
$ws.Range("F3").Value = 1.025503485680593
$ws.Range("I3").Value = 1.032482231056694
$ws.Range("J3").Value = 1.031177321845912
$ws.Range("K3").Value = 1.032874632585069
$ws.Range("L3").Value = 1.029214901078311
$ws.Range("M3").Value = 1.028138375758012
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.026934971342464
$ws.Range("D4").Value = 1.030660565935474
$ws.Range("E4").Value = 1.027062656855838
$ws.Range("F4").Value = 1.026479985932026
$ws.Range("I4").Value = 1.03259848946209
$ws.Range("J4").Value = 1.031509686634143
$ws.Range("K4").Value = 1.033158108933363
$ws.Range("L4").Value = 1.029569435292668
$ws.Range("M4").Value = 1.028988266592922
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.027173002000022
$ws.Range("D5").Value = 1.030832107136805
$ws.Range("E5").Value = 1.027264455693699
$ws.Range("F5").Value = 1.026890556546701
$ws.Range("I5").Value = 1.032647098525973
$ws.Range("J5").Value = 1.031649304348574
$ws.Range("K5").Value = 1.033277147422895
$ws.Range("L5").Value = 1.029718450032269
$ws.Range("M5").Value = 1.029345493173619
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.027212970816792
$ws.Range("D6").Value = 1.030860911356725
$ws.Range("E6").Value = 1.027298346040234
$ws.Range("F6").Value = 1.026959496206853
$ws.Range("I6").Value = 1.032655244596935
$ws.Range("J6").Value = 1.031672740407389
$ws.Range("K6").Value = 1.033297126579846
$ws.Range("L6").Value = 1.029743468402982
$ws.Range("M6").Value = 1.02940546922977
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.026938151757118
$ws.Range("D7").Value = 1.030662857963007
$ws.Range("E7").Value = 1.027065352803054
$ws.Range("F7").Value = 1.026485471795389
$ws.Range("I7").Value = 1.032599140024211
$ws.Range("J7").Value = 1.031511552639446
$ws.Range("K7").Value = 1.033159700062453
$ws.Range("L7").Value = 1.029571426559156
$ws.Range("M7").Value = 1.028993040132974
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.025789807018749
$ws.Range("D8").Value = 1.029835279393524
$ws.Range("E8").Value = 1.026092561238318
$ws.Range("F8").Value = 1.024504539659673
$ws.Range("I8").Value = 1.03236237560627
$ws.Range("J8").Value = 1.030836894762309
$ws.Range("K8").Value = 1.032584135590784
$ws.Range("L8").Value = 1.028852054364983
$ws.Range("M8").Value = 1.027268571499022
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.023767644504359
$ws.Range("D9").Value = 1.028377993833413
$ws.Range("E9").Value = 1.024382516626543
$ws.Range("F9").Value = 1.02101516219413
$ws.Range("I9").Value = 1.031936678954321
$ws.Range("J9").Value = 1.029644619602196
$ws.Range("K9").Value = 1.031565637707621
$ws.Range("L9").Value = 1.027583463769089
$ws.Range("M9").Value = 1.024227410010787
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.022420483115819
$ws.Range("D10").Value = 1.027407206774109
$ws.Range("E10").Value = 1.023245321255945
$ws.Range("F10").Value = 1.018689508897479
$ws.Range("I10").Value = 1.031647190476261
$ws.Range("J10").Value = 1.028847475053409
$ws.Range("K10").Value = 1.030883789958763
$ws.Range("L10").Value = 1.026737118447284
$ws.Range("M10").Value = 1.022198117753418
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.02183737702411
$ws.Range("D11").Value = 1.026987031218647
$ws.Range("E11").Value = 1.02275358587549
$ws.Range("F11").Value = 1.017682558812418
$ws.Range("I11").Value = 1.031520492088799
$ws.Range("J11").Value = 1.028501762057551
$ws.Range("K11").Value = 1.030587872124912
$ws.Range("L11").Value = 1.026370500899462
$ws.Range("M11").Value = 1.021318923492813
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.021620819386394
$ws.Range("D12").Value = 1.02683098746269
$ws.Range("E12").Value = 1.022571036080069
$ws.Range("F12").Value = 1.01730853860371
$ws.Range("I12").Value = 1.031473228355481
$ws.Range("J12").Value = 1.028373267330379
$ws.Range("K12").Value = 1.03047785430962
$ws.Range("L12").Value = 1.026234301544342
$ws.Range("M12").Value = 1.020992272568825
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.021667270194368
$ws.Range("D13").Value = 1.026864458090243
$ws.Range("E13").Value = 1.022610188974348
$ws.Range("F13").Value = 1.01738876702376
$ws.Range("I13").Value = 1.031483375738419
$ws.Range("J13").Value = 1.028400833571281
$ws.Range("K13").Value = 1.030501458057609
$ws.Range("L13").Value = 1.026263517724264
$ws.Range("M13").Value = 1.02106234394648
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.021819475608254
$ws.Range("D14").Value = 1.026974132013807
$ws.Range("E14").Value = 1.022738494154446
$ws.Range("F14").Value = 1.017651642094744
$ws.Range("I14").Value = 1.031516589379753
$ws.Range("J14").Value = 1.028491142305175
$ws.Range("K14").Value = 1.03057878007111
$ws.Range("L14").Value = 1.026359243049819
$ws.Range("M14").Value = 1.021291924062699
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.021913258892492
$ws.Range("D15").Value = 1.027041709494152
$ws.Range("E15").Value = 1.022817560832394
$ws.Range("F15").Value = 1.01781360872925
$ws.Range("I15").Value = 1.031537026598347
$ws.Range("J15").Value = 1.02854677370013
$ws.Range("K15").Value = 1.03062640737613
$ws.Range("L15").Value = 1.026418219775858
$ws.Range("M15").Value = 1.021433365348029
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.022459186686942
$ws.Range("D16").Value = 1.027435096342843
$ws.Range("E16").Value = 1.023277970460826
$ws.Range("F16").Value = 1.018756337964193
$ws.Range("I16").Value = 1.031655570648438
$ws.Range("J16").Value = 1.028870407441053
$ws.Range("K16").Value = 1.030903414899692
$ws.Range("L16").Value = 1.02676144663982
$ws.Range("M16").Value = 1.022256456135106
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.022801692882924
$ws.Range("D17").Value = 1.027681906763555
$ws.Range("E17").Value = 1.023566955268178
$ws.Range("F17").Value = 1.019347703092766
$ws.Range("I17").Value = 1.031729569455788
$ws.Range("J17").Value = 1.029073268812649
$ws.Range("K17").Value = 1.031076994429246
$ws.Range("L17").Value = 1.026976705415614
$ws.Range("M17").Value = 1.022772623053222
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.023001492512661
$ws.Range("D18").Value = 1.027825884575072
$ws.Range("E18").Value = 1.023735580553874
$ws.Range("F18").Value = 1.019692643766912
$ws.Range("I18").Value = 1.031772601635202
$ws.Range("J18").Value = 1.029191541901547
$ws.Range("K18").Value = 1.031178175446175
$ws.Range("L18").Value = 1.027102248253101
$ws.Range("M18").Value = 1.023073646789194
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.023069622670281
$ws.Range("D19").Value = 1.027874980224734
$ws.Range("E19").Value = 1.023793088476459
$ws.Range("F19").Value = 1.019810261032826
$ws.Range("I19").Value = 1.031787252427435
$ws.Range("J19").Value = 1.029231861053738
$ws.Range("K19").Value = 1.031212664523428
$ws.Range("L19").Value = 1.027145052734709
$ws.Range("M19").Value = 1.023176280145215
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.022764942977026
$ws.Range("D20").Value = 1.027655424528026
$ws.Range("E20").Value = 1.023535943169338
$ws.Range("F20").Value = 1.019284254481996
$ws.Range("I20").Value = 1.031721643536614
$ws.Range("J20").Value = 1.029051509137722
$ws.Range("K20").Value = 1.031058377707361
$ws.Range("L20").Value = 1.026953611627432
$ws.Range("M20").Value = 1.022717248191797
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.021774653968343
$ws.Range("D21").Value = 1.026941834992789
$ws.Range("E21").Value = 1.02270070863764
$ws.Range("F21").Value = 1.017574231816845
$ws.Range("I21").Value = 1.031506814373332
$ws.Range("J21").Value = 1.028464550882709
$ws.Range("K21").Value = 1.030556013443553
$ws.Range("L21").Value = 1.026331054911207
$ws.Range("M21").Value = 1.021224320656898
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.021152216416243
$ws.Range("D22").Value = 1.026493336240646
$ws.Range("E22").Value = 1.022176157652847
$ws.Range("F22").Value = 1.016499104223017
$ws.Range("I22").Value = 1.031370572344432
$ws.Range("J22").Value = 1.028095035345997
$ws.Range("K22").Value = 1.030239574340785
$ws.Range("L22").Value = 1.025939505916215
$ws.Range("M22").Value = 1.020285198477666
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.021482163611206
$ws.Range("D23").Value = 1.026731078148558
$ws.Range("E23").Value = 1.022454175474018
$ws.Range("F23").Value = 1.017069048421856
$ws.Range("I23").Value = 1.031442907694298
$ws.Range("J23").Value = 1.028290967150499
$ws.Range("K23").Value = 1.030407379804077
$ws.Range("L23").Value = 1.026147084953182
$ws.Range("M23").Value = 1.020783089901524
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.022781548610935
$ws.Range("D24").Value = 1.027667390656419
$ws.Range("E24").Value = 1.023549956000056
$ws.Range("F24").Value = 1.019312924154016
$ws.Range("I24").Value = 1.03172522532021
$ws.Range("J24").Value = 1.0290613415604
$ws.Range("K24").Value = 1.031066790003412
$ws.Range("L24").Value = 1.02696404675938
$ws.Range("M24").Value = 1.022742269854413
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.024290257077637
$ws.Range("D25").Value = 1.028754611257463
$ws.Range("E25").Value = 1.024824108352123
$ws.Range("F25").Value = 1.021917125289942
$ws.Range("I25").Value = 1.032047735924609
$ws.Range("J25").Value = 1.029953257143898
$ws.Range("K25").Value = 1.031829448066556
$ws.Range("L25").Value = 1.029214901078311
$ws.Range("M25").Value = 1.025013933357988
